$wb = $excel.ActiveWorkbook

# --- Sheet "altmed_retin_dermatitis": complete the per-indication ATC codes ---
# The "D07AB.." code cells in column B previously pointed at stray duplicate
# strings (with a trailing non-breaking space). Re-enter them as the clean,
# canonical code text so they match the codes used elsewhere in the workbook.
$wsDerm = $wb.Worksheets.Item("altmed_retin_dermatitis")

$wsDerm.Range("B2").Value = "D07AB01"
$wsDerm.Range("B3").Value = "D07AB02"
$wsDerm.Range("B4").Value = "D07AB03"
$wsDerm.Range("B5").Value = "D07AB04"
$wsDerm.Range("B6").Value = "D07AB05"
$wsDerm.Range("B7").Value = "D07AB06"
$wsDerm.Range("B8").Value = "D07AB07"
$wsDerm.Range("B9").Value = "D07AB08"
$wsDerm.Range("B10").Value = "D07AB09"
$wsDerm.Range("B11").Value = "D07AB10"
$wsDerm.Range("B12").Value = "D07AB11"
$wsDerm.Range("B13").Value = "D07AB19"
$wsDerm.Range("B14").Value = "D07AB21"
$wsDerm.Range("B15").Value = "D07AB30"

# --- Update selections / active sheet ---
# Move the selection on "altmed_retin_acne" off its old spot and make
# "altmed_retin_dermatitis" the active (selected) tab with its own selection.
$wsAcne = $wb.Worksheets.Item("altmed_retin_acne")
[void]$wsAcne.Activate()
[void]$wsAcne.Range("B2:B34").Select()

[void]$wsDerm.Activate()
[void]$wsDerm.Range("B2:B22").Select()

$wb.Save()
